# Re-run of the Case_4_144 (380 kV) line-loading scenario:
# loading_percent values for rows 2-25 (time steps 0-23) change in columns
# B, C, D, E, G, H, J, K, L. Columns A, F, I, M, N, O are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers for B, C, D, E, G, H, J, K, L
$cols = @(2, 3, 4, 5, 7, 8, 10, 11, 12)

# One row per worksheet row (2..25), values ordered to match $cols
$data = @(
    ,@(20.06401820495559, 5.220170146773357, 14.00620387218371, 14.3227799686162, 78.41004658587147, 26.44857315120652, 8.598654553216933, 15.56854269033155, 12.87185682475069)
    ,@(20.03116411070205, 5.155764915678615, 14.00015057389965, 14.3365572225787, 77.85908019617361, 26.3820023714264, 8.606929258942662, 15.56607233207755, 12.89227590668433)
    ,@(20.01673281822343, 5.11457694044593, 13.99883127220349, 14.34671526786602, 77.52751221480729, 26.3435023659408, 8.612309036474347, 15.56904645954179, 12.90681417990492)
    ,@(20.01230129764294, 5.097381607753773, 13.99889757932483, 14.35128220071737, 77.39416466880235, 26.32841519691191, 8.614576778870518, 15.57138883432355, 12.9132418460286)
    ,@(20.01165312204031, 5.094501585958774, 13.99894509154002, 14.35206636038472, 77.37213145370596, 26.32594646696812, 8.614957899021816, 15.57184607359074, 12.91433954793864)
    ,@(20.01666717843288, 5.114346698841618, 13.9988297199271, 14.34677512806838, 77.52570658062878, 26.34329645220723, 8.612339314278721, 15.56907347187017, 12.90689882817879)
    ,@(20.05150233275299, 5.198305495175726, 14.0036198292558, 14.32717797062871, 78.21872445527886, 26.4251296443253, 8.601445740020786, 15.56675972051093, 12.87848209865486)
    ,@(20.165056569129, 5.349824469855842, 14.03197865076289, 14.3022137976806, 79.62757849863993, 26.60422416111963, 8.582445887980725, 15.59777270133279, 12.8386288393905)
    ,@(20.27555782582798, 5.453040103052696, 14.06428424694069, 14.2920630414093, 80.68815447872998, 26.74683850189701, 8.569912170212465, 15.64206773308062, 12.81901692230326)
    ,@(20.33156922755872, 5.498216211413239, 14.08144613510168, 14.2892186443491, 81.17508344814594, 26.81403703292802, 8.564516661351545, 15.66683726829885, 12.81219138770457)
    ,@(20.35359155346084, 5.515066085423737, 14.08829678616504, 14.28839597737107, 81.36001825187272, 26.83981002304914, 8.562517305410514, 15.6768753452159, 12.8099077456417)
    ,@(20.34881277198336, 5.511448635668019, 14.08680577618294, 14.28856184629074, 81.32016653499645, 26.83424495887519, 8.562945957628978, 15.67468428623069, 12.81038618666692)
    ,@(20.33336482774849, 5.499607627841112, 14.08200270834068, 14.28914586694148, 81.19028777576499, 26.81615089656998, 8.564351296430001, 15.66764994847707, 12.81199748019094)
    ,@(20.32400782994848, 5.492321105128525, 14.07910641991502, 14.28953671538375, 81.1108014082708, 26.80511003030612, 8.565217805474038, 15.66342676648418, 12.8130236360468)
    ,@(20.27201176804506, 5.450051821284494, 14.06321206696778, 14.29228456636467, 80.65641415571928, 26.74249301392633, 8.570270921209634, 15.64054145780497, 12.81950513493067)
    ,@(20.24157549976389, 5.423664494065156, 14.05409109379291, 14.29442408447544, 80.37873790050854, 26.70466907641151, 8.573449099719047, 15.62768130586409, 12.82401793132755)
    ,@(20.22461091722558, 5.408320020782754, 14.04907723993496, 14.29582162041773, 80.21945227797187, 26.68313283990165, 8.575305932358553, 15.62071970236539, 12.82681088714782)
    ,@(20.21896041669308, 5.403095988715543, 14.0474196106845, 14.29632348999607, 80.16559707750048, 26.67587889476211, 8.575939581919949, 15.61843753733623, 12.82779043283091)
    ,@(20.244759540638, 5.426490785890666, 14.05503801515044, 14.29417905504244, 80.40825357910506, 26.70867287796816, 8.573107795138837, 15.62900528789575, 12.82351711756014)
    ,@(20.33788034271248, 5.50309261646827, 14.08340396122471, 14.28896742502625, 81.22842229368382, 26.82145676652836, 8.563937327190835, 15.66969828649539, 12.81151603740191)
    ,@(20.40346591368161, 5.551655646896477, 14.10399190301804, 14.2870441225659, 81.76758305803197, 26.89706604597581, 8.558199132649857, 15.70012780600769, 12.80542710184925)
    ,@(20.36803418429486, 5.525874510859159, 14.09281722818304, 14.28793514551193, 81.47956847530777, 26.8565408588749, 8.561238432127128, 15.68353832083286, 12.80851648922995)
    ,@(20.2433183725134, 5.425213561367409, 14.05460919565743, 14.29428931108119, 80.39490843023721, 26.70686210686488, 8.5732620065505, 15.62840536995877, 12.82374291717305)
    ,@(20.1295405110613, 5.310250877295521, 14.02228471423104, 14.30752735731267, 79.24160517044552, 26.55381196723566, 8.587334465791692, 15.58559041169683, 12.84771151602902)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Cells.Item($rowNum, $cols[$c]).Value = $rowVals[$c]
    }
}
